$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1882.78
$ws.Range("I15").Value = 1882.78
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 5648.34
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -5479.34
# Row 17
$ws.Range("H17").Value = 2836.875
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2836.875
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8510.625
$ws.Range("N17").Value = -8846.625
# Row 58
$ws.Range("H58").Value = 7299.5
$ws.Range("I58").Value = 400
$ws.Range("J58").Value = 10749.25
$ws.Range("K58").Value = 1200
$ws.Range("L58").Value = 32247.75
$ws.Range("M58").Value = -1050
$ws.Range("N58").Value = -32547.75
# Row 86
$ws.Range("H86").Value = 4977
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 4977
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 4977
$ws.Range("N86").Value = -7223
$ws.Range("M86").ClearContents()
# Row 89
$ws.Range("H89").Value = 4977
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 4977
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 24885
$ws.Range("N89").Value = -36117
$ws.Range("M89").ClearContents()
# Row 92
$ws.Range("H92").Value = 884.7143
$ws.Range("I92").Value = 803.1
$ws.Range("J92").Value = 1088.75
$ws.Range("K92").Value = 803.1
$ws.Range("L92").Value = 1088.75
$ws.Range("M92").Value = 444.9
$ws.Range("N92").Value = -3584.75
# Row 96
$ws.Range("H96").Value = 25000748
$ws.Range("I96").Value = 994.5
$ws.Range("J96").Value = 50000500
$ws.Range("K96").Value = 2983.5
$ws.Range("L96").Value = 150001500
$ws.Range("M96").Value = -1610.5
$ws.Range("N96").Value = -150004246
# Row 98
$ws.Range("H98").Value = 690.36365
$ws.Range("I98").Value = 768.75
$ws.Range("J98").Value = 481.33334
$ws.Range("K98").Value = 768.75
$ws.Range("L98").Value = 481.33334
$ws.Range("M98").Value = 729.25
$ws.Range("N98").Value = -3477.33334
# Row 103
$ws.Range("H103").Value = 1249.75
$ws.Range("I103").Value = 1500
$ws.Range("J103").Value = 999.5
$ws.Range("K103").Value = 4500
$ws.Range("L103").Value = 2998.5
$ws.Range("M103").Value = -3914
$ws.Range("N103").Value = -4170.5
# Row 104
$ws.Range("H104").Value = 1066.3334
$ws.Range("I104").Value = 1066.3334
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 3199.0002
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -1452.0002
# Row 122
$ws.Range("H122").Value = 690.36365
$ws.Range("I122").Value = 768.75
$ws.Range("J122").Value = 481.33334
$ws.Range("K122").Value = 2306.25
$ws.Range("L122").Value = 1444.00002
$ws.Range("M122").Value = 143.75
$ws.Range("N122").Value = -6344.000019999999
# Row 132
$ws.Range("H132").Value = 961.55554
$ws.Range("I132").Value = 961.55554
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2884.66662
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -354.66662
# Row 137
$ws.Range("H137").Value = 2467.9707
$ws.Range("I137").Value = 1488.3684
$ws.Range("J137").Value = 3708.8
$ws.Range("K137").Value = 4465.1052
$ws.Range("L137").Value = 11126.4
$ws.Range("M137").Value = -1915.1052
$ws.Range("N137").Value = -16226.4
# Row 138
$ws.Range("H138").Value = 4024.2273
$ws.Range("I138").Value = 3607.8333
$ws.Range("J138").Value = 4523.9
$ws.Range("K138").Value = 10823.4999
$ws.Range("L138").Value = 13571.7
$ws.Range("M138").Value = -5683.499899999999
$ws.Range("N138").Value = -23851.7

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6784.5063
$ws.Range("I32").Value = 6468.705
$ws.Range("J32").Value = 14995.333
$ws.Range("K32").Value = 6468.705
$ws.Range("L32").Value = 14995.333
$ws.Range("M32").Value = -6181.705
# Row 61
$ws.Range("H61").Value = 1307.4286
$ws.Range("I61").Value = 907.4643
$ws.Range("J61").Value = 2907.2856
$ws.Range("K61").Value = 907.4643
$ws.Range("L61").Value = 2907.2856
$ws.Range("M61").Value = -695.4643
# Row 132
$ws.Range("H132").Value = 2525.6667
$ws.Range("I132").Value = 1921.5714
$ws.Range("J132").Value = 3935.2222
$ws.Range("K132").Value = 5764.7142
$ws.Range("L132").Value = 11805.6666
$ws.Range("M132").Value = -3234.7142
$ws.Range("N132").Value = -16865.6666
# Row 133
$ws.Range("H133").Value = 75751.14
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 75751.14
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 75751.14
$ws.Range("N133").Value = -80811.14
# Row 136
$ws.Range("H136").Value = 1307.4286
$ws.Range("I136").Value = 907.4643
$ws.Range("J136").Value = 2907.2856
$ws.Range("K136").Value = 2722.3929
$ws.Range("L136").Value = 8721.856800000001
$ws.Range("M136").Value = -172.3928999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 37
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
# Row 86
$ws.Range("H86").Value = 4335.077
$ws.Range("I86").Value = 3564.3635
$ws.Range("J86").Value = 8574
$ws.Range("K86").Value = 3564.3635
$ws.Range("L86").Value = 8574
$ws.Range("M86").Value = -2441.3635
$ws.Range("N86").Value = -10820
# Row 89
$ws.Range("H89").Value = 4335.077
$ws.Range("I89").Value = 3564.3635
$ws.Range("J89").Value = 8574
$ws.Range("K89").Value = 17821.8175
$ws.Range("L89").Value = 42870
$ws.Range("M89").Value = -12205.8175
$ws.Range("N89").Value = -54102
# Row 94
$ws.Range("H94").Value = 808.5
$ws.Range("I94").Value = 528
$ws.Range("J94").Value = 3333
$ws.Range("K94").Value = 528
$ws.Range("L94").Value = 3333
$ws.Range("M94").Value = -77
# Row 134
$ws.Range("H134").Value = 2620.1035
$ws.Range("I134").Value = 2218.5908
$ws.Range("J134").Value = 3882
$ws.Range("K134").Value = 6655.7724
$ws.Range("L134").Value = 11646
$ws.Range("M134").Value = -4120.7724
$ws.Range("N134").Value = -16716

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2508.9285
$ws.Range("I31").Value = 2634
$ws.Range("J31").Value = 2342.1667
$ws.Range("K31").Value = 2634
$ws.Range("L31").Value = 2342.1667
$ws.Range("M31").Value = -2339
# Row 34
$ws.Range("H34").Value = 2508.9285
$ws.Range("I34").Value = 2634
$ws.Range("J34").Value = 2342.1667
$ws.Range("K34").Value = 2634
$ws.Range("L34").Value = 2342.1667
$ws.Range("M34").Value = -2432
# Row 58
$ws.Range("H58").Value = 2925.25
$ws.Range("I58").Value = 2862.5557
$ws.Range("J58").Value = 3005.8572
$ws.Range("K58").Value = 2862.5557
$ws.Range("L58").Value = 3005.8572
$ws.Range("M58").Value = -2659.5557
# Row 59
$ws.Range("H59").Value = 60127
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 60127
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 60127
$ws.Range("N59").Value = -62417
# Row 68
$ws.Range("H68").Value = 70147.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 70147.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 70147.5
$ws.Range("N68").Value = -71645.5
# Row 71
$ws.Range("H71").Value = 70147.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 70147.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 210442.5
$ws.Range("N71").Value = -217930.5
# Row 136
$ws.Range("H136").Value = 2925.25
$ws.Range("I136").Value = 2862.5557
$ws.Range("J136").Value = 3005.8572
$ws.Range("K136").Value = 8587.667099999999
$ws.Range("L136").Value = 9017.571599999999
$ws.Range("M136").Value = -6037.667099999999
# Row 141
$ws.Range("H141").Value = 51045.6
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 51045.6
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 51045.6
$ws.Range("N141").Value = -61405.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 25700138
$ws.Range("I4").Value = 49848484
$ws.Range("J4").Value = 16040800
$ws.Range("K4").Value = 149545452
$ws.Range("L4").Value = 48122400
$ws.Range("M4").Value = -149545340

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4399.625
$ws.Range("I80").Value = 3518.8
$ws.Range("J80").Value = 5867.6665
$ws.Range("K80").Value = 3518.8
$ws.Range("L80").Value = 5867.6665
$ws.Range("M80").Value = -2520.8
$ws.Range("N80").Value = -7863.6665
# Row 83
$ws.Range("H83").Value = 4399.625
$ws.Range("I83").Value = 3518.8
$ws.Range("J83").Value = 5867.6665
$ws.Range("K83").Value = 17594
$ws.Range("L83").Value = 29338.3325
$ws.Range("M83").Value = -12602
$ws.Range("N83").Value = -39322.3325
# Row 126
$ws.Range("H126").Value = 2150
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -13340

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 1840.5
$ws.Range("I82").Value = 1770.25
$ws.Range("J82").Value = 1981
$ws.Range("K82").Value = 1770.25
$ws.Range("L82").Value = 1981
$ws.Range("M82").Value = -1409.25
# Row 85
$ws.Range("H85").Value = 1840.5
$ws.Range("I85").Value = 1770.25
$ws.Range("J85").Value = 1981
$ws.Range("K85").Value = 1770.25
$ws.Range("L85").Value = 1981
$ws.Range("M85").Value = -522.25
# Row 93
$ws.Range("H93").Value = 1912.375
$ws.Range("I93").Value = 1912.375
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1912.375
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -664.375
# Row 100
$ws.Range("H100").Value = 1600
$ws.Range("I100").Value = 200
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 200
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = 341
$ws.Range("N100").Value = -4082
# Row 132
$ws.Range("H132").Value = 5366.4443
$ws.Range("I132").Value = 5083.5
$ws.Range("J132").Value = 5932.3335
$ws.Range("K132").Value = 15250.5
$ws.Range("L132").Value = 17797.0005
$ws.Range("M132").Value = -12720.5
$ws.Range("N132").Value = -22857.0005
# Row 136
$ws.Range("H136").Value = 3546.276
$ws.Range("I136").Value = 3991.5715
$ws.Range("J136").Value = 2377.375
$ws.Range("K136").Value = 11974.7145
$ws.Range("L136").Value = 7132.125
$ws.Range("M136").Value = -9424.7145
$ws.Range("N136").Value = -12232.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("H94").Value = 35000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 35000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 35000
$ws.Range("N94").Value = -36802
# Row 100
$ws.Range("H100").Value = 1419.6
$ws.Range("I100").Value = 1419.6
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2839.2
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2298.2
# Row 122
$ws.Range("H122").Value = 4499.25
$ws.Range("I122").Value = 4499.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13497.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11047.75
# Row 132
$ws.Range("H132").Value = 2120
$ws.Range("I132").Value = 1160.4445
$ws.Range("J132").Value = 4998.6665
$ws.Range("K132").Value = 3481.3335
$ws.Range("L132").Value = 14995.9995
$ws.Range("M132").Value = -951.3335000000002
$ws.Range("N132").Value = -20055.9995
# Row 136
$ws.Range("H136").Value = 3592.3635
$ws.Range("I136").Value = 3426.6155
$ws.Range("J136").Value = 3831.7778
$ws.Range("K136").Value = 10279.8465
$ws.Range("L136").Value = 11495.3334
$ws.Range("M136").Value = -7729.8465
